$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at row 81 (shifts old rows 81-83 down to 82-84),
# reusing the blank-row formatting that was previously on row 81.
$ws1.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new regression-test summary
# line ("Demo-Baseline 2010-18 C340").
$ws1.Range("A81").Value = "CW3M"
$ws1.Range("B81").Value = "Demo-Baseline 2010-18 C340"
$ws1.Range("C81").Value = "2010-18"
$ws1.Range("D81").Value = 1211.2894424444446
$ws1.Range("E81").Value = 1890.2624783333331
$ws1.Range("F81").Value = 1.0618724444444443
$ws1.Range("G81").Value = 270.41205844444437
$ws1.Range("H81").Value = 9.775355222222224
$ws1.Range("I81").Value = 6.0645103333333337
$ws1.Range("J81").Value = 8.145128999999999
$ws1.Range("K81").Value = 693.57212322222222
$ws1.Range("L81").Value = 82.308506444444433
$ws1.Range("M81").Value = 1421.1092122222226
$ws1.Range("N81").Value = 1183.8681913333335
$ws1.Range("O81").Value = 5611.0926921111113
$ws1.Range("P81").Value = 27227.338324777778
$ws1.Range("Q81").Value = 0.13744488888888892
$ws1.Range("R81").Value = 0.000021777777777777772

# Move the selection to B82 (the now-blank row that used to be row 81),
# and make this sheet the active tab.
$ws1.Range("B82").Select()
